# Apply the per-battery segregation edit to the "Analysis Results" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Starting/Ending SoC (%) values (rows 6 & 7) ---
$ws.Range("B6").Value = 96
$ws.Range("B7").Value = 32

# --- Relabel rows 8-10 with units ---
$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("A10").Value = "Total SOC consumed(%)"

# --- Relabel rows 12-15 with units ---
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"
$ws.Range("A15").Value = "Regenerative Effectiveness(kWh)"

# --- Rows 16 & 17: swap Lowest/Highest Cell Voltage labels + values ---
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("B16").Value = 3.334
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("B17").Value = 3.09

$ws.Range("A18").Value = "Difference in Cell Voltage(V)"

# --- Rows 19-21: temperature labels + fill in the missing difference value ---
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("B21").Value = 10

# --- Rows 22-27: BMS/motor temperature labels ---
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"

# --- Rows 28 & 29: swap lowest/highest cell temp labels ---
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"

$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"

# --- Row 31: was "Maximum BMS Temperature in C" -> now "Battery Voltage(V)" ---
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("B31").Value = 55

# --- Row 32: was "Battery Voltage" -> now "Total energy charged(kWh)" ---
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("B32").Value = 1.386001252777778

# --- Row 33: was "Total energy charged in kWh" -> now "Electricity consumption units(kW)" ---
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("B33").Value = 0.0000001082373764390855

# --- Row 34: was "Electricity consumption units in kW" -> now "Idling time percentage" ---
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("B34").Value = 3.994024208055123

# --- Row 35: was "Idling time percentage" -> now "Time spent in 0-10 km/h" ---
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("B35").Value = 23.91536327327053

# --- Row 36: was "Time spent in 0-10 km/h" -> now "Time spent in 10-20 km/h" ---
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("B36").Value = 4.378182261654319

# --- Row 37: was "Time spent in 10-20 km/h" -> now "Time spent in 20-30 km/h" ---
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("B37").Value = 8.856977346870332

# --- Row 38: was "Time spent in 20-30 km/h" -> now "Time spent in 30-40 km/h" ---
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("B38").Value = 18.68044757462118

# --- Row 39: was "Time spent in 30-40 km/h" -> now "Time spent in 40-50 km/h" ---
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("B39").Value = 14.79313393701027

# --- Row 40: was "Time spent in 40-50 km/h" -> now "Time spent in 50-60 km/h" ---
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("B40").Value = 19.07984999542669

# --- Row 41: was "Time spent in 50-60 km/h" -> now "Time spent in 60-70 km/h" ---
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("B41").Value = 5.884325741638465

# --- Row 42: was "Time spent in 60-70 km/h" -> now "Time spent in 70-80 km/h" ---
$ws.Range("A42").Value = "Time spent in 70-80 km/h"
$ws.Range("B42").Value = 0

# --- Row 43 (new): "Time spent in 80-90 km/h" ---
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
